$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 534.4545000000001
$ws.Range("I18").Value = 410.2
$ws.Range("J18").Value = 1777
$ws.Range("K18").Value = 410.2
$ws.Range("L18").Value = 1777
$ws.Range("M18").Value = -126.2
$ws.Range("N18").Value = -2345

$ws.Range("H113").Value = 2081.5881
$ws.Range("I113").Value = 1629
$ws.Range("J113").Value = 2398.4
$ws.Range("K113").Value = 1629
$ws.Range("L113").Value = 2398.4
$ws.Range("M113").Value = 1625
$ws.Range("N113").Value = -8906.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 895.7619
$ws.Range("I2").Value = 741.4
$ws.Range("J2").Value = 1281.6666
$ws.Range("K2").Value = 741.4
$ws.Range("L2").Value = 1281.6666
$ws.Range("M2").Value = -628.4
$ws.Range("N2").Value = -1507.6666

$ws.Range("H23").Value = 85006.25
$ws.Range("I23").Value = 80006
$ws.Range("J23").Value = 100007
$ws.Range("K23").Value = 80006
$ws.Range("L23").Value = 100007
$ws.Range("M23").Value = -79747
$ws.Range("N23").Value = -100525

$ws.Range("H28").Value = 5400
$ws.Range("I28").Value = 5400
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 5400
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -5208
$ws.Range("N28").ClearContents()

$ws.Range("H32").Value = 397193.38
$ws.Range("I32").Value = 425962.7
$ws.Range("J32").Value = 16000
$ws.Range("K32").Value = 425962.7
$ws.Range("L32").Value = 16000
$ws.Range("M32").Value = -425675.7
$ws.Range("N32").Value = -16574

$ws.Range("H99").Value = 5400
$ws.Range("I99").Value = 5400
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 5400
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -2405
$ws.Range("N99").Value = -2405

$ws.Range("H116").Value = 895.7619
$ws.Range("I116").Value = 741.4
$ws.Range("J116").Value = 1281.6666
$ws.Range("K116").Value = 741.4
$ws.Range("L116").Value = 1281.6666
$ws.Range("M116").Value = 1552.6
$ws.Range("N116").Value = -5869.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 895.7619
$ws.Range("I3").Value = 741.4
$ws.Range("J3").Value = 1281.6666
$ws.Range("K3").Value = 741.4
$ws.Range("L3").Value = 1281.6666
$ws.Range("M3").Value = -627.4
$ws.Range("N3").Value = -1509.6666

$ws.Range("H86").Value = 3194.2856
$ws.Range("I86").Value = 2800
$ws.Range("J86").Value = 3260
$ws.Range("K86").Value = 2800
$ws.Range("L86").Value = 3260
$ws.Range("M86").Value = -1677
$ws.Range("N86").Value = -5506

$ws.Range("H89").Value = 3194.2856
$ws.Range("I89").Value = 2800
$ws.Range("J89").Value = 3260
$ws.Range("K89").Value = 14000
$ws.Range("L89").Value = 16300
$ws.Range("M89").Value = -8384
$ws.Range("N89").Value = -27532

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 206.11111
$ws.Range("I19").Value = 206.11111
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 206.11111
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -36.11111
$ws.Range("N19").ClearContents()

$ws.Range("H22").Value = 603.8
$ws.Range("I22").Value = 284.15384
$ws.Range("J22").Value = 2681.5
$ws.Range("K22").Value = 284.15384
$ws.Range("L22").Value = 2681.5
$ws.Range("M22").Value = 65.84616
$ws.Range("N22").Value = -3381.5

$ws.Range("H24").Value = 206.11111
$ws.Range("I24").Value = 206.11111
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 206.11111
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -36.11111
$ws.Range("N24").ClearContents()

$ws.Range("H31").Value = 2152.73
$ws.Range("I31").Value = 1272.7046
$ws.Range("J31").Value = 2844.1785
$ws.Range("K31").Value = 1272.7046
$ws.Range("L31").Value = 2844.1785
$ws.Range("M31").Value = -977.7046
$ws.Range("N31").Value = -3434.1785

$ws.Range("H34").Value = 2152.73
$ws.Range("I34").Value = 1272.7046
$ws.Range("J34").Value = 2844.1785
$ws.Range("K34").Value = 1272.7046
$ws.Range("L34").Value = 2844.1785
$ws.Range("M34").Value = -1070.7046
$ws.Range("N34").Value = -3248.1785

$ws.Range("H132").Value = 31255484
$ws.Range("I132").Value = 55563470
$ws.Range("J132").Value = 2354.6428
$ws.Range("K132").Value = 166690410
$ws.Range("L132").Value = 7063.928400000001
$ws.Range("M132").Value = -166687880
$ws.Range("N132").Value = -12123.9284

$ws.Range("H134").Value = 2372.6155
$ws.Range("I134").Value = 2340
$ws.Range("J134").Value = 2461.1428
$ws.Range("K134").Value = 7020
$ws.Range("L134").Value = 7383.428400000001
$ws.Range("M134").Value = -4485
$ws.Range("N134").Value = -12453.4284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 873.4666999999999
$ws.Range("I5").Value = 961.3
$ws.Range("J5").Value = 697.8
$ws.Range("K5").Value = 2883.9
$ws.Range("L5").Value = 2093.4
$ws.Range("M5").Value = -2771.9
$ws.Range("N5").Value = -2317.4

$ws.Range("H132").Value = 48385.477
$ws.Range("I132").Value = 679.6429000000001
$ws.Range("J132").Value = 143797.14
$ws.Range("K132").Value = 6116.7861
$ws.Range("L132").Value = 1294174.26
$ws.Range("M132").Value = -3586.7861
$ws.Range("N132").Value = -1299234.26

$ws.Range("H135").Value = 873.4666999999999
$ws.Range("I135").Value = 961.3
$ws.Range("J135").Value = 697.8
$ws.Range("K135").Value = 8651.699999999999
$ws.Range("L135").Value = 6280.2
$ws.Range("M135").Value = -6116.699999999999
$ws.Range("N135").Value = -11350.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 36338.426
$ws.Range("I132").Value = 63683.445
$ws.Range("J132").Value = 3524.4
$ws.Range("K132").Value = 191050.335
$ws.Range("L132").Value = 10573.2
$ws.Range("M132").Value = -188520.335
$ws.Range("N132").Value = -15633.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()

$ws.Range("H107").Value = 441.88235
$ws.Range("I107").Value = 429.66666
$ws.Range("J107").Value = 471.2
$ws.Range("K107").Value = 1288.99998
$ws.Range("L107").Value = 1413.6
$ws.Range("M107").Value = 631.0000199999999
$ws.Range("N107").Value = -5253.6
